$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to text
# so Excel does not silently convert "28.30" -> 28.3, "6.90" -> 6.9, etc.
$textCells = @('D5', 'D6', 'D14', 'D18', 'D19', 'D20', 'D21', 'D23', 'D28', 'D35', 'D37', 'D38', 'D39', 'D41', 'D43', 'D44', 'D45', 'D48', 'D49', 'D50')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.857.75'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '3.450.53'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '580.24'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').Value = '148.31'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('E10').Value = '  -2.32%  '
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').Value = '4.044.01'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').Value = '28.30'
$ws.Range('E14').Value = '  -4.51%  '
$ws.Range('D15').Value = '3.453.10'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '62.963.14'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '6.47'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').Value = '14.61'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '9.13'
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').Value = '388.64'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '74.84'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '3.594.13'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').Value = '7.65'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -2.76%  '
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('E33').Value = '  -6.08%  '
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('D35').Value = '1.63'
$ws.Range('E35').Value = '  +3.95%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '7.02'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').Value = '31.80'
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').Value = '170.30'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = '3.486.99'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').Value = '0.0781'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '42.75'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').Value = '1.72'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('D45').Value = '4.35'
$ws.Range('E45').Value = '  -3.26%  '
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('D47').Value = '2.567.67'
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = '6.90'
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '2.27'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').Value = '22.70'
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('E51').Value = '  +0.12%  '

# Restore default (General) style on the forced-text cells so no stray
# cell style/number-format artifacts are left behind in the output.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
